# Adds a "2022" column (column S) of data to the sheet, mirroring the
# formatting of the existing "2021" column (column R), then moves the
# active selection to T15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new value for column S (year 2022 data).
# Row 34 is intentionally omitted: it stays an empty, styled cell (like R34).
# (Values are written in plain decimal - not scientific-notation - form,
# since they parse to the identical IEEE-754 double either way.)
$values = @{
    3  = 2022
    4  = 0.15686557910355481
    5  = 0.18747863920572591
    6  = 0.12556588018347117
    7  = 0.051313356512815066
    8  = 0.060745743331368028
    9  = 0.042060988433228183
    10 = 0.092022006630303563
    11 = 0.078942235953699605
    12 = 0.10098382728705417
    13 = 0.097010038673425045
    14 = 0.12657756598786343
    15 = 0.067310604785784003
    16 = 0.12618253497302423
    17 = 0.15767275020694549
    18 = 0.083781780685077176
    19 = 0.089790167285988584
    20 = 0.11543537913568107
    21 = 0.064489306438090949
    22 = 0.077235413540471365
    23 = 0.1335826876836762
    24 = 0.021874179718260566
    25 = 0.13849188927432132
    26 = 0.15541703258327452
    27 = 0.12135301021830269
    28 = 0.4304881257025327
    29 = 0.49554896622979544
    30 = 0.35193780867878632
    31 = 0.21076296192215821
    32 = 0.25905990040586052
    33 = 0.1647039446594746
    35 = 0
    36 = 0.1
    37 = 0.2
}

for ($row = 3; $row -le 37; $row++) {
    # Column R (18) already carries the formatting we want column S (19) to have.
    $srcCell = $ws.Cells.Item($row, 18)
    $dstCell = $ws.Cells.Item($row, 19)

    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    if ($values.ContainsKey($row)) {
        $dstCell.Value = $values[$row]
    }
}

$excel.CutCopyMode = 0

# Move the selection, matching the saved workbook state.
$ws.Range("T15").Select() | Out-Null
